# Update "想去人数" (F column) counts on the 展览 and 全部类型 sheets.
# Mapping of worksheet row number -> new F-column value.
$updates = @{
    2  = 145
    7  = 1323
    9  = 344
    10 = 433
    15 = 116
    16 = 279
    17 = 319
    19 = 1765
    20 = 73
    26 = 4251
    28 = 287
    29 = 1115
    32 = 619
    33 = 24
    34 = 312
    36 = 157
}

$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
